$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 14-17: summary labels + aggregate stats ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Format B14 (bold, size 12, vertically centered), then copy that format
# down to B15:B17 so they all share one style/font entry.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

$ws.Range("A14:B17").RowHeight = 15.6

# --- Row 12: average of column J ("|S*|/n") ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$avgCell = $ws.Range("J12")
$avgCell.Font.Bold = $true
$avgCell.Font.Name = "Calibri"

# --- Sheet-level bookkeeping to mirror the saved workbook state ---
$ws.Range("A14:B17").Select()
